# basic module added (#8)
# Rework the inventory sheet: add category/suppliers/customers columns,
# replace the sample rows with the new fixture data, and drop the
# now-removed rows 4-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -----------------------------------------------
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "category"
$ws.Range("D1").Value = "suppliers"
$ws.Range("E1").Value = "customers"
$ws.Range("F1").Value = "quantity"
$ws.Range("G1").Value = "price"
$ws.Range("H1").Value = "created_at"
$ws.Range("I1").Value = "updated_at"

# --- Row 2 --------------------------------------------------------------
$ws.Range("A2").Value = "5fbe7df6-340c-460c-b064-3dbf08aa63d7"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "cee738ea-c625-4c80-add1-b875a50aa499"
$ws.Range("D2").Value = '["2f0e4354-8965-4cc4-8b57-b685d58039c0"]'
$ws.Range("E2").Value = '["e664bdf5-f208-4757-9be2-425c3dc3b6d6"]'
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = "2024-11-17T12:26:27.664712"
$ws.Range("I2").Value = 0

# --- Row 3 --------------------------------------------------------------
$ws.Range("A3").Value = "0c291ffb-a8d8-47fb-8cf9-235fe6a54b62"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "cee738ea-c625-4c80-add1-b875a50aa499"
$ws.Range("D3").Value = '["ee77754b-f795-469d-a40d-998282b919e6", "2f0e4354-8965-4cc4-8b57-b685d58039c0"]'
$ws.Range("E3").Value = '["e664bdf5-f208-4757-9be2-425c3dc3b6d6"]'
$ws.Range("F3").Value = 123
$ws.Range("G3").Value = 1234
$ws.Range("H3").Value = "2024-11-17T12:32:05.218225"
$ws.Range("I3").Value = 0

# --- Drop the old rows 4-7 (data no longer present in the fixture) -----
$ws.Range("A4:I7").EntireRow.Delete()
